$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PastrySalesData")

$ws.Range("A15").Value = "SconE"
$ws.Range("A16").Value = "SCONE"
